$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 202.3573813333333
$ws.Cells.Item(2, 8).Value = 607.072144
$ws.Cells.Item(2, 9).Value = 0.4567501787232752
$ws.Cells.Item(2, 10).Value = 0.4567501787232752
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.065175
$ws.Cells.Item(2, 14).Value = 0.195525
$ws.Cells.Item(2, 15).Value = 0.009404016458916581
$ws.Cells.Item(2, 16).Value = 0.009404016458916581
$ws.Cells.Item(2, 17).Value = 13.1886423284
$ws.Cells.Item(2, 18).Value = 118.6977809556
$ws.Cells.Item(2, 19).Value = 0.004295286198326771
$ws.Cells.Item(2, 20).Value = 0.00429528619832677

$ws.Cells.Item(3, 7).Value = 202.3573813333333
$ws.Cells.Item(3, 8).Value = 607.072144
$ws.Cells.Item(3, 9).Value = 0.4567501787232752
$ws.Cells.Item(3, 10).Value = 0.4567501787232752
$ws.Cells.Item(3, 13).Value = 6.718514333333332
$ws.Cells.Item(3, 14).Value = 20.155543
$ws.Cells.Item(3, 15).Value = 0.969405744075698
$ws.Cells.Item(3, 16).Value = 0.969405744075698
$ws.Cells.Item(3, 17).Value = 1359.540966943799
$ws.Cells.Item(3, 18).Value = 12235.86870249419
$ws.Cells.Item(3, 19).Value = 0.4427762468619447
$ws.Cells.Item(3, 20).Value = 0.4427762468619446

$ws.Cells.Item(4, 7).Value = 202.3573813333333
$ws.Cells.Item(4, 8).Value = 607.072144
$ws.Cells.Item(4, 9).Value = 0.4567501787232752
$ws.Cells.Item(4, 10).Value = 0.4567501787232752
$ws.Cells.Item(4, 11).Value = 1
$ws.Cells.Item(4, 12).Value = 0.3333333333333333
$ws.Cells.Item(4, 13).Value = 0.14686
$ws.Cells.Item(4, 14).Value = 0.44058
$ws.Cells.Item(4, 15).Value = 0.02119023946538534
$ws.Cells.Item(4, 16).Value = 0.02119023946538533
$ws.Cells.Item(4, 17).Value = 29.71820502261334
$ws.Cells.Item(4, 18).Value = 267.46384520352
$ws.Cells.Item(4, 19).Value = 0.009678645663003752
$ws.Cells.Item(4, 20).Value = 0.00967864566300375

$ws.Cells.Item(5, 9).Value = 0.5138029191566978
$ws.Cells.Item(5, 10).Value = 0.5138029191566978
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 0.065175
$ws.Cells.Item(5, 14).Value = 0.195525
$ws.Cells.Item(5, 15).Value = 0.009404016458916581
$ws.Cells.Item(5, 16).Value = 0.009404016458916581
$ws.Cells.Item(5, 17).Value = 14.83603782485
$ws.Cells.Item(5, 18).Value = 133.52434042365
$ws.Cells.Item(5, 19).Value = 0.004831811108388972
$ws.Cells.Item(5, 20).Value = 0.004831811108388972

$ws.Cells.Item(6, 9).Value = 0.5138029191566978
$ws.Cells.Item(6, 10).Value = 0.5138029191566978
$ws.Cells.Item(6, 13).Value = 6.718514333333332
$ws.Cells.Item(6, 14).Value = 20.155543
$ws.Cells.Item(6, 15).Value = 0.969405744075698
$ws.Cells.Item(6, 16).Value = 0.969405744075698
$ws.Cells.Item(6, 17).Value = 1529.361454179213
$ws.Cells.Item(6, 18).Value = 13764.25308761292
$ws.Cells.Item(6, 19).Value = 0.4980835011533644
$ws.Cells.Item(6, 20).Value = 0.4980835011533644

$ws.Cells.Item(7, 9).Value = 0.5138029191566978
$ws.Cells.Item(7, 10).Value = 0.5138029191566978
$ws.Cells.Item(7, 11).Value = 1
$ws.Cells.Item(7, 12).Value = 0.3333333333333333
$ws.Cells.Item(7, 13).Value = 0.14686
$ws.Cells.Item(7, 14).Value = 0.44058
$ws.Cells.Item(7, 15).Value = 0.02119023946538534
$ws.Cells.Item(7, 16).Value = 0.02119023946538533
$ws.Cells.Item(7, 17).Value = 33.43031093145333
$ws.Cells.Item(7, 18).Value = 300.87279838308
$ws.Cells.Item(7, 19).Value = 0.01088760689494445
$ws.Cells.Item(7, 20).Value = 0.01088760689494445

$ws.Cells.Item(8, 5).Value = 2
$ws.Cells.Item(8, 6).Value = 0.6666666666666666
$ws.Cells.Item(8, 7).Value = 0.1759946666666667
$ws.Cells.Item(8, 8).Value = 0.527984
$ws.Cells.Item(8, 9).Value = 0.0003972456795234369
$ws.Cells.Item(8, 10).Value = 0.0003972456795234369
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 0.065175
$ws.Cells.Item(8, 14).Value = 0.195525
$ws.Cells.Item(8, 15).Value = 0.009404016458916581
$ws.Cells.Item(8, 16).Value = 0.009404016458916581
$ws.Cells.Item(8, 17).Value = 0.0114704524
$ws.Cells.Item(8, 18).Value = 0.1032340716
$ws.Cells.Item(8, 19).Value = 0.000003735704908471902
$ws.Cells.Item(8, 20).Value = 0.000003735704908471902

$ws.Cells.Item(9, 5).Value = 2
$ws.Cells.Item(9, 6).Value = 0.6666666666666666
$ws.Cells.Item(9, 7).Value = 0.1759946666666667
$ws.Cells.Item(9, 8).Value = 0.527984
$ws.Cells.Item(9, 9).Value = 0.0003972456795234369
$ws.Cells.Item(9, 10).Value = 0.0003972456795234369
$ws.Cells.Item(9, 13).Value = 6.718514333333332
$ws.Cells.Item(9, 14).Value = 20.155543
$ws.Cells.Item(9, 15).Value = 0.969405744075698
$ws.Cells.Item(9, 16).Value = 0.969405744075698
$ws.Cells.Item(9, 17).Value = 1.182422690590222
$ws.Cells.Item(9, 18).Value = 10.641804215312
$ws.Cells.Item(9, 19).Value = 0.0003850922435392736
$ws.Cells.Item(9, 20).Value = 0.0003850922435392736

$ws.Cells.Item(10, 5).Value = 2
$ws.Cells.Item(10, 6).Value = 0.6666666666666666
$ws.Cells.Item(10, 7).Value = 0.1759946666666667
$ws.Cells.Item(10, 8).Value = 0.527984
$ws.Cells.Item(10, 9).Value = 0.0003972456795234369
$ws.Cells.Item(10, 10).Value = 0.0003972456795234369
$ws.Cells.Item(10, 11).Value = 1
$ws.Cells.Item(10, 12).Value = 0.3333333333333333
$ws.Cells.Item(10, 13).Value = 0.14686
$ws.Cells.Item(10, 14).Value = 0.44058
$ws.Cells.Item(10, 15).Value = 0.02119023946538534
$ws.Cells.Item(10, 16).Value = 0.02119023946538533
$ws.Cells.Item(10, 17).Value = 0.02584657674666667
$ws.Cells.Item(10, 18).Value = 0.23261919072
$ws.Cells.Item(10, 19).Value = 0.000008417731075691348
$ws.Cells.Item(10, 20).Value = 0.000008417731075691348

$ws.Cells.Item(11, 7).Value = 12.07540333333333
$ws.Cells.Item(11, 8).Value = 36.22621
$ws.Cells.Item(11, 9).Value = 0.02725594981667764
$ws.Cells.Item(11, 10).Value = 0.02725594981667764
$ws.Cells.Item(11, 11).Value = 1
$ws.Cells.Item(11, 12).Value = 0.3333333333333333
$ws.Cells.Item(11, 13).Value = 0.065175
$ws.Cells.Item(11, 14).Value = 0.195525
$ws.Cells.Item(11, 15).Value = 0.009404016458916581
$ws.Cells.Item(11, 16).Value = 0.009404016458916581
$ws.Cells.Item(11, 17).Value = 0.78701441225
$ws.Cells.Item(11, 18).Value = 7.083129710250001
$ws.Cells.Item(11, 19).Value = 0.0002563154006794409
$ws.Cells.Item(11, 20).Value = 0.0002563154006794409

$ws.Cells.Item(12, 7).Value = 12.07540333333333
$ws.Cells.Item(12, 8).Value = 36.22621
$ws.Cells.Item(12, 9).Value = 0.02725594981667764
$ws.Cells.Item(12, 10).Value = 0.02725594981667764
$ws.Cells.Item(12, 13).Value = 6.718514333333332
$ws.Cells.Item(12, 14).Value = 20.155543
$ws.Cells.Item(12, 15).Value = 0.969405744075698
$ws.Cells.Item(12, 16).Value = 0.969405744075698
$ws.Cells.Item(12, 17).Value = 81.1287703757811
$ws.Cells.Item(12, 18).Value = 730.15893338203
$ws.Cells.Item(12, 19).Value = 0.02642207431252627
$ws.Cells.Item(12, 20).Value = 0.02642207431252627

$ws.Cells.Item(13, 7).Value = 12.07540333333333
$ws.Cells.Item(13, 8).Value = 36.22621
$ws.Cells.Item(13, 9).Value = 0.02725594981667764
$ws.Cells.Item(13, 10).Value = 0.02725594981667764
$ws.Cells.Item(13, 11).Value = 1
$ws.Cells.Item(13, 12).Value = 0.3333333333333333
$ws.Cells.Item(13, 13).Value = 0.14686
$ws.Cells.Item(13, 14).Value = 0.44058
$ws.Cells.Item(13, 15).Value = 0.02119023946538534
$ws.Cells.Item(13, 16).Value = 0.02119023946538533
$ws.Cells.Item(13, 17).Value = 1.773393733533334
$ws.Cells.Item(13, 18).Value = 15.9605436018
$ws.Cells.Item(13, 19).Value = 0.0005775601034719247
$ws.Cells.Item(13, 20).Value = 0.0005775601034719246

$ws.Cells.Item(14, 7).Value = 0.7692486666666666
$ws.Cells.Item(14, 8).Value = 2.307746
$ws.Cells.Item(14, 9).Value = 0.001736306645537542
$ws.Cells.Item(14, 10).Value = 0.001736306645537542
$ws.Cells.Item(14, 11).Value = 1
$ws.Cells.Item(14, 12).Value = 0.3333333333333333
$ws.Cells.Item(14, 13).Value = 0.065175
$ws.Cells.Item(14, 14).Value = 0.195525
$ws.Cells.Item(14, 15).Value = 0.009404016458916581
$ws.Cells.Item(14, 16).Value = 0.009404016458916581
$ws.Cells.Item(14, 17).Value = 0.05013578184999999
$ws.Cells.Item(14, 18).Value = 0.45122203665
$ws.Cells.Item(14, 19).Value = 0.00001632825627236128
$ws.Cells.Item(14, 20).Value = 0.00001632825627236128

$ws.Cells.Item(15, 7).Value = 0.7692486666666666
$ws.Cells.Item(15, 8).Value = 2.307746
$ws.Cells.Item(15, 9).Value = 0.001736306645537542
$ws.Cells.Item(15, 10).Value = 0.001736306645537542
$ws.Cells.Item(15, 13).Value = 6.718514333333332
$ws.Cells.Item(15, 14).Value = 20.155543
$ws.Cells.Item(15, 15).Value = 0.969405744075698
$ws.Cells.Item(15, 16).Value = 0.969405744075698
$ws.Cells.Item(15, 17).Value = 5.168208192897554
$ws.Cells.Item(15, 18).Value = 46.51387373607799
$ws.Cells.Item(15, 19).Value = 0.0016831856356609
$ws.Cells.Item(15, 20).Value = 0.0016831856356609

$ws.Cells.Item(16, 7).Value = 0.7692486666666666
$ws.Cells.Item(16, 8).Value = 2.307746
$ws.Cells.Item(16, 9).Value = 0.001736306645537542
$ws.Cells.Item(16, 10).Value = 0.001736306645537542
$ws.Cells.Item(16, 11).Value = 1
$ws.Cells.Item(16, 12).Value = 0.3333333333333333
$ws.Cells.Item(16, 13).Value = 0.14686
$ws.Cells.Item(16, 14).Value = 0.44058
$ws.Cells.Item(16, 15).Value = 0.02119023946538534
$ws.Cells.Item(16, 16).Value = 0.02119023946538533
$ws.Cells.Item(16, 17).Value = 0.1129718591866667
$ws.Cells.Item(16, 18).Value = 1.01674673268
$ws.Cells.Item(16, 19).Value = 0.00003679275360428044
$ws.Cells.Item(16, 20).Value = 0.00003679275360428044

$ws.Cells.Item(17, 5).Value = 1
$ws.Cells.Item(17, 6).Value = 0.3333333333333333
$ws.Cells.Item(17, 7).Value = 0.02543033333333333
$ws.Cells.Item(17, 8).Value = 0.076291
$ws.Cells.Item(17, 9).Value = 0.00005739997828821048
$ws.Cells.Item(17, 10).Value = 0.00005739997828821048
$ws.Cells.Item(17, 11).Value = 1
$ws.Cells.Item(17, 12).Value = 0.3333333333333333
$ws.Cells.Item(17, 13).Value = 0.065175
$ws.Cells.Item(17, 14).Value = 0.195525
$ws.Cells.Item(17, 15).Value = 0.009404016458916581
$ws.Cells.Item(17, 16).Value = 0.009404016458916581
$ws.Cells.Item(17, 17).Value = 0.001657421975
$ws.Cells.Item(17, 18).Value = 0.014916797775
$ws.Cells.Item(17, 19).Value = 0.0000005397903405637858
$ws.Cells.Item(17, 20).Value = 0.0000005397903405637858

$ws.Cells.Item(18, 5).Value = 1
$ws.Cells.Item(18, 6).Value = 0.3333333333333333
$ws.Cells.Item(18, 7).Value = 0.02543033333333333
$ws.Cells.Item(18, 8).Value = 0.076291
$ws.Cells.Item(18, 9).Value = 0.00005739997828821048
$ws.Cells.Item(18, 10).Value = 0.00005739997828821048
$ws.Cells.Item(18, 13).Value = 6.718514333333332
$ws.Cells.Item(18, 14).Value = 20.155543
$ws.Cells.Item(18, 15).Value = 0.969405744075698
$ws.Cells.Item(18, 16).Value = 0.969405744075698
$ws.Cells.Item(18, 17).Value = 0.1708540590014444
$ws.Cells.Item(18, 18).Value = 1.537686531013
$ws.Cells.Item(18, 19).Value = 0.00005564386866241159
$ws.Cells.Item(18, 20).Value = 0.00005564386866241159

$ws.Cells.Item(19, 5).Value = 1
$ws.Cells.Item(19, 6).Value = 0.3333333333333333
$ws.Cells.Item(19, 7).Value = 0.02543033333333333
$ws.Cells.Item(19, 8).Value = 0.076291
$ws.Cells.Item(19, 9).Value = 0.00005739997828821048
$ws.Cells.Item(19, 10).Value = 0.00005739997828821048
$ws.Cells.Item(19, 11).Value = 1
$ws.Cells.Item(19, 12).Value = 0.3333333333333333
$ws.Cells.Item(19, 13).Value = 0.14686
$ws.Cells.Item(19, 14).Value = 0.44058
$ws.Cells.Item(19, 15).Value = 0.02119023946538534
$ws.Cells.Item(19, 16).Value = 0.02119023946538533
$ws.Cells.Item(19, 17).Value = 0.003734698753333333
$ws.Cells.Item(19, 18).Value = 0.03361228878
$ws.Cells.Item(19, 19).Value = 0.000001216319285235099
$ws.Cells.Item(19, 20).Value = 0.000001216319285235099

